$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.178.70'
$ws.Range("E2").Value = '  +0.37%  '
$ws.Range("D3").Value = '2.022.08'
$ws.Range("E3").Value = '  +0.09%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = "'228.50"
$ws.Range("E5").Value = '  +0.65%  '
$ws.Range("D6").Value = "'0.611"
$ws.Range("E6").Value = '  +0.26%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = "'55.93"
$ws.Range("E8").Value = '  +1.71%  '
$ws.Range("E9").Value = '  -1.04%  '
$ws.Range("D10").Value = "'0.0780"
$ws.Range("E10").Value = '  -1.78%  '
$ws.Range("D11").Value = "'0.103"
$ws.Range("E11").Value = '  -1.88%  '
$ws.Range("D12").Value = '2.321.04'
$ws.Range("E12").Value = '  +0.08%  '
$ws.Range("D13").Value = "'14.25"
$ws.Range("E13").Value = '  -0.30%  '
$ws.Range("D14").Value = "'20.14"
$ws.Range("E14").Value = '  -2.42%  '
$ws.Range("D15").Value = "'0.738"
$ws.Range("E15").Value = '  -0.97%  '
$ws.Range("D16").Value = "'5.17"
$ws.Range("E16").Value = '  +0.24%  '
$ws.Range("D17").Value = '2.022.62'
$ws.Range("E17").Value = '  -0.68%  '
$ws.Range("D18").Value = '37.119.23'
$ws.Range("E18").Value = '  +0.50%  '
$ws.Range("D19").Value = "'6.16"
$ws.Range("D20").Value = "'68.67"
$ws.Range("E20").Value = '  -0.31%  '
$ws.Range("D21").Value = '0.0₃0816'
$ws.Range("E21").Value = '  -1.79%  '
$ws.Range("D22").Value = "'222.66"
$ws.Range("E22").Value = '  -1.53%  '
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("E24").Value = '  +1.86%  '
$ws.Range("D25").Value = "'2.22"
$ws.Range("E25").Value = '  -2.41%  '
$ws.Range("D26").Value = "'163.73"
$ws.Range("E26").Value = '  -2.12%  '
$ws.Range("D27").Value = "'9.01"
$ws.Range("E27").Value = '  -3.42%  '
$ws.Range("E28").Value = '  +2.65%  '
$ws.Range("D29").Value = "'18.66"
$ws.Range("E29").Value = '  -0.70%  '
$ws.Range("E30").Value = '  -2.13%  '
$ws.Range("E31").Value = '  -0.10%  '
$ws.Range("D32").Value = "'4.44"
$ws.Range("E32").Value = '  -0.96%  '
$ws.Range("D33").Value = "'0.0605"
$ws.Range("E33").Value = '  -1.32%  '
$ws.Range("B34").Value = 'WEMIXToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D34").Value = "'1.99"
$ws.Range("E34").Value = '  +9.01%  '
$ws.Range("B35").Value = 'InternetComputer(DFINITY)'
$ws.Range("C35").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D35").Value = "'4.42"
$ws.Range("E35").Value = '  -0.67%  '
$ws.Range("E36").Value = '  -2.56%  '
$ws.Range("D37").Value = "'3.20"
$ws.Range("E37").Value = '  +0.71%  '
$ws.Range("E38").Value = '  +0.04%  '
$ws.Range("D39").Value = "'5.57"
$ws.Range("E39").Value = '  +2.32%  '
$ws.Range("D40").Value = "'4.41"
$ws.Range("E40").Value = '  +19.37%  '
$ws.Range("D41").Value = '1.467.73'
$ws.Range("E41").Value = '  -2.24%  '
$ws.Range("E42").Value = '  -3.06%  '
$ws.Range("D43").Value = "'2.82"
$ws.Range("E43").Value = '  -0.07%  '
$ws.Range("D44").Value = "'93.60"
$ws.Range("E44").Value = '  -1.67%  '
$ws.Range("E45").Value = '  -1.75%  '
$ws.Range("D46").Value = "'16.13"
$ws.Range("E46").Value = '  -4.79%  '
$ws.Range("D47").Value = "'1.10"
$ws.Range("E47").Value = '  -2.90%  '
$ws.Range("E48").Value = '  +0.22%  '
$ws.Range("D49").Value = "'7.10"
$ws.Range("E49").Value = '  -2.09%  '
$ws.Range("E50").Value = '  +0.54%  '
$ws.Range("D51").Value = '2.209.73'
$ws.Range("E51").Value = '  +0.03%  '
